$wb = $excel.ActiveWorkbook

# --- Step 1: duplicate the "2022-Q2" sheet so the old data survives on its own tab ---
# Copy sheet #2 ("2022-Q2") to just after itself; the copy becomes sheet #3.
$wsQ2 = $wb.Worksheets.Item(2)
$wsQ2.Copy($null, $wsQ2)

# The original sheet (position 2) will hold the NEW "2022-Q3" figures; the copy
# (position 3) keeps the original "2022-Q2" figures. Rename the original first so the
# copy (still carrying the old default "2022-Q2 (2)" style name) can reclaim "2022-Q2".
$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Name = "2022-Q3"
$wsQ2old = $wb.Worksheets.Item(3)
$wsQ2old.Name = "2022-Q2"

# --- Step 2: overwrite the renamed "2022-Q3" sheet's data row with the new quarter's figures ---
# (Fund code / rank stay the same; name + numeric-looking figures are stored as text,
# matching the source data which keeps them as inline strings.)
$wsQ3.Range("C2").Value = "华安国际龙头（DAX）ETF（QDII）"

$wsQ3.Range("D2:G2").NumberFormat = "@"
$wsQ3.Range("D2").Value = "5.54"
$wsQ3.Range("E2").Value = "93.57"
$wsQ3.Range("F2").Value = "9.77"
$wsQ3.Range("G2").Value = "0.5413"
$wsQ3.Range("D2:G2").Style = "Normal"

# Match the header-row formatting used on the "总计" sheet (bold, bordered, centered)
# and the same style for the leading A2 marker cell.
$wsTotal = $wb.Worksheets.Item(1)
$wsTotal.Range("B1:D1").Copy()
$wsQ3.Range("B1:H1").PasteSpecial(-4122)
$wsTotal.Range("A2").Copy()
$wsQ3.Range("A2").PasteSpecial(-4122)

# --- Step 3: update the "总计" (totals) sheet ---
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("D2").Value = 0.54

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q2"
$wsTotal.Range("C3").Value = 1
$wsTotal.Range("D3").Value = 0.57

# Carry the same styling used in row 2 down into the new row 3.
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)
